# "Ajuste de la Queratina / % a Marinela / Menor % a los profesionales de Tocador"
#
# The payroll sheet gains three new columns (Porc_trans, Cost_trans, Porc_producto)
# inserted before the old "Valor_producto" column, plus a new "Valor_Neto" column
# inserted before "Part_profesional". The old layout was:
#   A Fecha de Pago | B Nombre cliente | C Servicio/Producto | D Prestador/Vendedor
#   E Precio | F Valor_producto | G Part_profesional | H Revisar
# The new layout is:
#   A..E (unchanged) | F Porc_trans | G Cost_trans | H Porc_producto | I Valor_producto
#   J Valor_Neto | K Part_profesional | L Revisar
#
# Row 2/3 also swap their service (Balayage now leads, Blower moves down) and every
# data row gets recalculated Porc_trans/Cost_trans/Porc_producto/Valor_producto/
# Valor_Neto/Part_profesional figures per the new commission rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("F1").Value = "Porc_trans"
$ws.Range("G1").Value = "Cost_trans"
$ws.Range("H1").Value = "Porc_producto"
$ws.Range("I1").Value = "Valor_producto"
$ws.Range("J1").Value = "Valor_Neto"
$ws.Range("K1").Value = "Part_profesional"
$ws.Range("L1").Value = "Revisar"

# ---- Row 2: now "Balayage Previa Valoracion & Diagnostico desde" ----
$ws.Range("C2").Value = "Balayage Previa Valoracion & Diagnostico desde"
$ws.Range("E2").Value = 365000
$ws.Range("F2").Value = 0.036
$ws.Range("G2").Value = 13140
$ws.Range("H2").Value = 0.1096767123287671
$ws.Range("I2").Value = 40032
$ws.Range("J2").Value = 311828
$ws.Range("K2").Value = 160718

# ---- Row 3: now "Blower  Cabello medio" ----
$ws.Range("C3").Value = "Blower  Cabello medio"
$ws.Range("E3").Value = 35000
$ws.Range("F3").Value = 0.036
$ws.Range("G3").Value = 1260
$ws.Range("H3").Value = 0.1166
$ws.Range("I3").Value = 4081
$ws.Range("J3").Value = 29659
$ws.Range("K3").Value = 19250

# ---- Row 4: Corte caballero ----
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.14575
$ws.Range("I4").Value = 4081
$ws.Range("J4").Value = 23919
$ws.Range("K4").Value = 15400

# ---- Row 5: Blower  Cabello medio ----
$ws.Range("F5").Value = 0.036
$ws.Range("G5").Value = 1260
$ws.Range("H5").Value = 0.1166
$ws.Range("I5").Value = 4081
$ws.Range("J5").Value = 29659
$ws.Range("K5").Value = 19250

# ---- Row 6: Blower  Cabello medio ----
$ws.Range("F6").Value = 0.036
$ws.Range("G6").Value = 1260
$ws.Range("H6").Value = 0.1166
$ws.Range("I6").Value = 4081
$ws.Range("J6").Value = 29659
$ws.Range("K6").Value = 19250

# ---- Row 7: Shampoo Dirigido ----
$ws.Range("F7").Value = 0.036
$ws.Range("G7").Value = 791.9999999999999
$ws.Range("H7").Value = 0.25
$ws.Range("I7").Value = 5500
$ws.Range("J7").Value = 15708
$ws.Range("K7").Value = 8800

# ---- Row 8: Base Global Tinte desde ----
$ws.Range("F8").Value = 0.036
$ws.Range("G8").Value = 5940
$ws.Range("H8").Value = 0.2222545454545455
$ws.Range("I8").Value = 36672
$ws.Range("J8").Value = 122388
$ws.Range("K8").Value = 54078

# ---- Row 9: Blower cabello extra largo ----
$ws.Range("F9").Value = 0.036
$ws.Range("G9").Value = 1980
$ws.Range("H9").Value = 0.0742
$ws.Range("I9").Value = 4081
$ws.Range("J9").Value = 48939
$ws.Range("K9").Value = 30250

# ---- Row 10: Queratina caballero - Desde ----
$ws.Range("F10").Value = 0.036
$ws.Range("G10").Value = 10800
$ws.Range("H10").Value = 0.2016
$ws.Range("I10").Value = 60480
$ws.Range("J10").Value = 228720
$ws.Range("K10").Value = 120736

# ---- Row 11: Corte caballero ----
$ws.Range("F11").Value = 0.036
$ws.Range("G11").Value = 1008
$ws.Range("H11").Value = 0.14575
$ws.Range("I11").Value = 4081
$ws.Range("J11").Value = 22911
$ws.Range("K11").Value = 15400

# ---- Row 12: Maquillaje Halloween ----
$ws.Range("F12").Value = 0.036
$ws.Range("G12").Value = 7199.999999999999
$ws.Range("H12").Value = 0.020405
$ws.Range("I12").Value = 4081
$ws.Range("J12").Value = 188719
$ws.Range("K12").Value = 110000

# ---- Row 13: Blower  Cabello medio ----
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0.1166
$ws.Range("I13").Value = 4081
$ws.Range("J13").Value = 30919
$ws.Range("K13").Value = 19250

# ---- Row 14: Corte Cabello Dama ----
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0.106
$ws.Range("I14").Value = 4081
$ws.Range("J14").Value = 34419
$ws.Range("K14").Value = 21175

# ---- Row 15: Ondas tubo o plancha ----
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0.2267222222222222
$ws.Range("I15").Value = 4081
$ws.Range("J15").Value = 13919
$ws.Range("K15").Value = 9900

# ---- Row 16: Blower  Cabello medio ----
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.1166
$ws.Range("I16").Value = 4081
$ws.Range("J16").Value = 30919
$ws.Range("K16").Value = 19250

# ---- Rows 17-19: fund deductions -- their single value moves from column G to column K ----
$ws.Range("K17").Value = -50000
$ws.Range("K18").Value = -17000
$ws.Range("K19").Value = -8000
$ws.Range("G17").ClearContents()
$ws.Range("G18").ClearContents()
$ws.Range("G19").ClearContents()
